$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38

# Row 3
$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 4
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48
$ws.Range("S3").Value = 1.57
$ws.Range("AC3").Value = 6
$ws.Range("AH3").Value = 19

# Row 4
$ws.Range("S4").Value = 1.62

# Row 5
$ws.Range("S5").Value = 1.5

# Row 6
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 5.25
$ws.Range("K6").Value = 1.91
$ws.Range("S6").Value = 1.62
$ws.Range("T6").Value = 2.2
$ws.Range("U6").Value = 2.5
$ws.Range("V6").Value = 1.5
$ws.Range("W6").Value = 4.75
$ws.Range("AL6").Value = 67
$ws.Range("AR6").Value = 81
$ws.Range("AT6").Value = 2.2

# Row 8
$ws.Range("G8").Value = 1.3
$ws.Range("H8").Value = 4.33
$ws.Range("J8").Value = 1.8
$ws.Range("K8").Value = 2.5
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 15
$ws.Range("O8").Value = 1.2
$ws.Range("P8").Value = 4.33
$ws.Range("Q8").Value = 1.67
$ws.Range("R8").Value = 2.15
$ws.Range("S8").Value = 1.3
$ws.Range("T8").Value = 3.4
$ws.Range("U8").Value = 1.95
$ws.Range("V8").Value = 1.8
$ws.Range("W8").Value = 7.5
$ws.Range("AA8").Value = 11
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 13
$ws.Range("AD8").Value = 9
$ws.Range("AE8").Value = 19
$ws.Range("AF8").Value = 51
$ws.Range("AG8").Value = 23
$ws.Range("AM8").Value = 301
$ws.Range("AP8").Value = 17
$ws.Range("AQ8").Value = 17
$ws.Range("AT8").Value = 3.4
$ws.Range("AU8").Value = 9
$ws.Range("AV8").Value = 51
$ws.Range("BB8").Value = 301

# Row 10
$ws.Range("M10").Value = 1.03
$ws.Range("N10").Value = 15

Write-Output "Done applying odds updates"
